$wb = $excel.ActiveWorkbook

# Update the zh-cn handback report: the first file (39b07019-...) has been
# handed back again, so refresh its Correspond Handoff/Handback datetimes.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 06:50:39"
$wsZhCn.Range("H2").Value = "2016-03-13 06:50:57"

# Update the de-de handback report the same way.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 06:50:43"
$wsDeDe.Range("H2").Value = "2016-03-13 06:51:03"
